$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")
$ws.Range("AS1").EntireColumn.Delete()
